$d = $word.ActiveDocument

# --- 1. Drop the stray spell-check bracket markers around "RobotContra" ---
# These <w:proofErr w:type="spellStart"/> / <w:proofErr w:type="spellEnd"/>
# markers are left over from a prior proofing pass and are cleared once the
# paragraph's content is regenerated. Rebuild the first paragraph in place
# (delete its range, reinsert an empty paragraph, retype the text) so the
# stale markers are not carried forward.
$p1 = $d.Paragraphs(1)
$p1.Range.Delete()
$d.Range(0, 0).InsertParagraphBefore()
$d.Paragraphs(1).Range.Text = "RobotContra"

# --- 2. Add the new "network function" paragraphs after "add child to parent" ---
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "search child"

$last2 = $d.Paragraphs.Last
$last2.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "handleinput"
